$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Cells.Item(14, 1).Value = " Originea musculus extensor digitorum: "
$ws.Cells.Item(14, 2).Value = "Epicondylus lateralis humeri"
$ws.Cells.Item(14, 3).Value = "Epicondylus medialis humeri"
$ws.Cells.Item(14, 4).Value = "Olecranon"
$ws.Cells.Item(14, 5).Value = "Tuberositas radii"
$ws.Cells.Item(14, 6).Value = "Fascia antebrachii"
$ws.Cells.Item(14, 9).Value = "Mușchiul extensor al degetelor își are originea pe epicondilul lateral al humerusului și pe fascia antebrațului. Distal mușchiul formează patru tendoane, care se inseră pe fața dorsală a falangelor într-un mod cu totul deosebit.Enunțuri corecte sunt „A” și „E”."

# Row 15
$ws.Cells.Item(15, 1).Value = "Muşchii mâinii, eminenţa tenară:"
$ws.Cells.Item(15, 2).Value = "Se disting muşchii palmari, dorsali şi interosoşi."
$ws.Cells.Item(15, 3).Value = "Sunt situaţi în trei planuri."
$ws.Cells.Item(15, 4).Value = "Toţi, cu excepţia adductorului policelului, au originea pe rândul proximal al oaselor carpiene."
$ws.Cells.Item(15, 5).Value = "M.opozant al policelului se inseră pe primul os metacarpian."
$ws.Cells.Item(15, 9).Value = "С. Самой глубокой является м-ца противопостаиляющая большой палец.проксимального ряда.Mușchii mâinii sunt amplasați numai pe fața palmară și în spațiile intermetacarpiene; pe fața dorsală a mâinii se află numai tendoanele mușchilor posteriori ai antebrațului.Mușchii mâinii formează trei grupuri:- lateral (eminența tenară);- medial (eminența hipotenară);- mediu (lombricalii, interosoșii palmari și dorsali).Mușchii eminenței tenare sunt dispuși în trei planuri. În primul plan se află abductorul scurt al policelui, în planul al doilea – apozantul policelui și flexorul scurt al policelui, în planul al treilea – adductorul policelui. Ei au originea unii pe oasele carpiene din rândul I, alții pe carpienele din rândul II, pe metacarpiene și pe retinaculul flexorilor și inserția sau pe osul metacarpian I (mușchiul apozant al policelui), sau pe falanga proximală a policelui. Astfel corecte sunt numai enunțurile „B” și „E”."

# Row 16
$ws.Cells.Item(16, 1).Value = "Muşchii mâinii, eminenţa tenară:"
$ws.Cells.Item(16, 2).Value = "Se disting muşchii palmari, dorsali şi interosoşi."
$ws.Cells.Item(16, 3).Value = "Sunt situaţi în trei planuri."
$ws.Cells.Item(16, 4).Value = "Toţi, cu excepţia adductorului policelului, au originea pe rândul proximal al oaselor carpiene."
$ws.Cells.Item(16, 5).Value = "M.opozant al policelului se inseră pe primul os metacarpian."

# Row 17
$ws.Cells.Item(17, 1).Value = "Muşchii mâinii, eminenţa tenară:"
$ws.Cells.Item(17, 2).Value = "Se disting muşchii palmari, dorsali şi interosoşi."
$ws.Cells.Item(17, 3).Value = "Sunt situaţi în trei planuri."
$ws.Cells.Item(17, 4).Value = "Toţi, cu excepţia adductorului policelului, au originea pe rândul proximal al oaselor carpiene."
$ws.Cells.Item(17, 5).Value = "M.opozant al policelului se inseră pe primul os metacarpian."
$ws.Cells.Item(17, 9).Value = "Mușchii mâinii sunt amplasați numai pe fața palmară și în spațiile intermetacarpiene; pe fața dorsală a mâinii se află numai tendoanele mușchilor posteriori ai antebrațului.Mușchii mâinii formează trei grupuri:- lateral (eminența tenară);- medial (eminența hipotenară);- mediu (lombricalii, interosoșii palmari și dorsali).Mușchii eminenței tenare sunt dispuși în trei planuri. În primul plan se află abductorul scurt al policelui, în planul al doilea – apozantul policelui și flexorul scurt al policelui, în planul al treilea – adductorul policelui. Ei au originea unii pe oasele carpiene din rândul I, alții pe carpienele din rândul II, pe metacarpiene și pe retinaculul flexorilor și inserția sau pe osul metacarpian I (mușchiul apozant al policelui), sau pe falanga proximală a policelui. Astfel corecte sunt numai enunțurile „B” și „E”."
